$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.001.84"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "3.076.02"
$ws.Range("E3").Value = "  -3.19%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -8.36%  "
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "3.075.09"
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "90.093.95"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("E16").Value = "  -4.83%  "
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.095.25"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000214"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.42%  "
$ws.Range("E23").Value = "  +5.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -15.02%  "
$ws.Range("D28").Value = "3.240.48"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.195"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.152"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "498.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.70%  "
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +52.09%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("E46").Value = "  -5.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.677"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "149.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("E51").Value = "  -6.30%  "
